$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("航天动力", "中国西电", "航天动力"),
    @("顺灏股份", "雪人集团", "永辉超市"),
    @("雪人集团", "摩尔线程-U", "航天发展"),
    @("摩尔线程-U", "华菱线缆", "雪人集团"),
    @("航天发展", "顺灏股份", "龙洲股份"),
    @("中国西电", "航天动力", "摩尔线程"),
    @("永辉超市", "龙洲股份", "平潭发展"),
    @("龙洲股份", "中超控股", "顺灏股份"),
    @("国机重装", "安泰科技", "中超控股"),
    @("平潭发展", "永辉超市", "再升科技"),
    @("西部材料", "国机重装", "博纳影业"),
    @("华菱线缆", "航天发展", "中国西电"),
    @("中超控股", "西部材料", "西部材料"),
    @("永鼎股份", "再升科技", "永鼎股份"),
    @("东方电气", "中能电气", "东百集团"),
    @("再升科技", "芯原股份", "华菱线缆"),
    @("安泰科技", "东方电气", "兆新股份"),
    @("芯原股份", "永鼎股份", "国机重装"),
    @("博纳影业", "博纳影业", "华映科技"),
    @("中能电气", "航天机电", "特变电工")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

Write-Output "Updated A2:C21 with new stock ranking order"
